$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($targetCell, $scratchCell, $text)
    $ws.Range($scratchCell).Formula = '="' + $text + '"'
    $ws.Range($scratchCell).Copy()
    $ws.Range($targetCell).PasteSpecial(-4163)
    $ws.Range($scratchCell).ClearContents()
}

Set-TextValue "P2" "D2" "320017960315"
Set-TextValue "P3" "D3" "320017959458"
Set-TextValue "P4" "D4" "320017959480"
Set-TextValue "P5" "D5" "320017959506"
Set-TextValue "P6" "D6" "320017959540"
Set-TextValue "P7" "D7" "320017959561"
Set-TextValue "P8" "D8" "320017959594"
Set-TextValue "P9" "D9" "320017959610"
Set-TextValue "P10" "D10" "320017959642"
Set-TextValue "P11" "D11" "320017959664"
Set-TextValue "P12" "D12" "320017959701"
Set-TextValue "P13" "D13" "320017959723"
Set-TextValue "P14" "D14" "320017959756"
Set-TextValue "P15" "D15" "320017959778"
Set-TextValue "P16" "D16" "320017959804"
Set-TextValue "P17" "D17" "320017959826"
Set-TextValue "Q18" "D18" "$85.66"
$ws.Range("R18").Value = "FAIL"
Set-TextValue "P18" "D18" "320017959860"
Set-TextValue "P19" "D19" "320017959881"
Set-TextValue "Q20" "D20" "$85.66"
$ws.Range("R20").Value = "FAIL"
Set-TextValue "P20" "D20" "320017959918"
Set-TextValue "P21" "D21" "320017959930"
Set-TextValue "Q22" "D22" "$233.07"
$ws.Range("R22").Value = "FAIL"
Set-TextValue "P22" "D22" "320017959962"
Set-TextValue "Q23" "D23" "$476.72"
$ws.Range("R23").Value = "FAIL"
Set-TextValue "P23" "D23" "320017959973"
Set-TextValue "Q24" "D24" "$306.68"
$ws.Range("R24").Value = "FAIL"
Set-TextValue "P24" "D24" "320017959984"
Set-TextValue "P25" "D25" "320017959995"
Set-TextValue "P26" "D26" "320017960006"
